# Generate Report for Handback
#
# Refreshes the localization-status report after a successful handback:
#   - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#     (this text is shared by the Overview sheet's per-locale status columns
#     and by the "Status" column on each per-locale detail sheet).
#   - The zh-cn and de-de detail sheets get fresh "Latest Handback DateTime"
#     stamps and their stale "Error Detail" (version-mismatch warning) is
#     cleared now that the handback is in sync.
#   - A couple of columns are widened on the detail/overview sheets to fit
#     the new text, and the now-unused "Error Detail" column is narrowed.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: per-locale status columns (zh-cn = E, de-de = F)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

# Widen to fit the longer status text.
$wsOverview.Columns("E:F").ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# zh-cn detail sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $statusText
$wsZhCn.Range("C3").Value = $statusText
$wsZhCn.Range("K2").Value = "2016-08-13 10:44:30"
$wsZhCn.Range("K3").Value = "2016-08-13 10:44:30"
$wsZhCn.Range("P2").Value = ""
$wsZhCn.Range("P3").Value = ""

$wsZhCn.Columns("C:C").ColumnWidth = 29.166666666666668
$wsZhCn.Columns("P:P").ColumnWidth = 12.833333333333334

# ---------------------------------------------------------------------
# de-de detail sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $statusText
$wsDeDe.Range("C3").Value = $statusText
$wsDeDe.Range("K2").Value = "2016-08-13 10:44:40"
$wsDeDe.Range("K3").Value = "2016-08-13 10:44:40"
$wsDeDe.Range("P2").Value = ""
$wsDeDe.Range("P3").Value = ""

$wsDeDe.Columns("C:C").ColumnWidth = 29.166666666666668
$wsDeDe.Columns("P:P").ColumnWidth = 12.833333333333334
